$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("input")
$wsOptions = $wb.Worksheets.Item("options")
$wsExplanation = $wb.Worksheets.Item("explanation")

# --- Update the "input" sheet's data row (row 2) ---
$wsInput.Range("A2").Value = "Nepal"
$wsInput.Range("D2").Value = "Network based on OSM online"
$wsInput.Range("E2").Value = "npl_admbnda_adm0_nd_20201117.shp"
$wsInput.Range("G2").ClearContents()
$wsInput.Range("H2").ClearContents()
$wsInput.Range("L2").Value = "drive"
$wsInput.Range("M2").Value = "motorway, trunk, primary, secondary"

# --- Update sheet selections / which tab is active ---
# "explanation" was the active tab before; its own in-sheet selection moves to A2.
$wsExplanation.Activate()
$wsExplanation.Range("A2").Select()

# "input" becomes the active (selected) tab, with the whole second row selected
# (active cell A2, selection spans A2:XFD2 - a full-row selection).
$wsInput.Activate()
$wsInput.Range("A2:XFD2").Select()
